$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay as literal text (matches original inlineStr cells),
# so force text format before assignment, then restore default "Normal" style so no
# extra style index lingers on the cell (keeps XML output clean/matching).

$d2 = $ws.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "26.307.22"
$d2.Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

$d3 = $ws.Range("D3")
$d3.NumberFormat = "@"
$d3.Value = "1.680.60"
$d3.Style = "Normal"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  +0.00%  "

$d5 = $ws.Range("D5")
$d5.NumberFormat = "@"
$d5.Value = "218.74"
$d5.Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

$d6 = $ws.Range("D6")
$d6.NumberFormat = "@"
$d6.Value = "0.5273"
$d6.Style = "Normal"
$ws.Range("E6").Value = "  +3.15%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  +2.13%  "

$d9 = $ws.Range("D9")
$d9.NumberFormat = "@"
$d9.Value = "0.06439"
$d9.Style = "Normal"
$ws.Range("E9").Value = "  +1.13%  "

$d10 = $ws.Range("D10")
$d10.NumberFormat = "@"
$d10.Value = "22.08"
$d10.Style = "Normal"
$ws.Range("E10").Value = "  +2.97%  "

$d11 = $ws.Range("D11")
$d11.NumberFormat = "@"
$d11.Value = "0.07513"
$d11.Style = "Normal"
$ws.Range("E11").Value = "  +1.74%  "

$d12 = $ws.Range("D12")
$d12.NumberFormat = "@"
$d12.Value = "1.710.48"
$d12.Style = "Normal"
$ws.Range("E12").Value = "  +2.36%  "

$d13 = $ws.Range("D13")
$d13.NumberFormat = "@"
$d13.Value = "4.550"
$d13.Style = "Normal"
$ws.Range("E13").Value = "  +0.29%  "

$d14 = $ws.Range("D14")
$d14.NumberFormat = "@"
$d14.Value = "0.5815"
$d14.Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "

$d15 = $ws.Range("D15")
$d15.NumberFormat = "@"
$d15.Value = "0.000008518"
$d15.Style = "Normal"
$ws.Range("E15").Value = "  -1.24%  "

$d16 = $ws.Range("D16")
$d16.NumberFormat = "@"
$d16.Value = "64.53"
$d16.Style = "Normal"
$ws.Range("E16").Value = "  +0.27%  "

$d17 = $ws.Range("D17")
$d17.NumberFormat = "@"
$d17.Value = "26.348.28"
$d17.Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "

$d18 = $ws.Range("D18")
$d18.NumberFormat = "@"
$d18.Value = "4.936"
$d18.Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "

$d19 = $ws.Range("D19")
$d19.NumberFormat = "@"
$d19.Value = "1.007"
$d19.Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "

$d20 = $ws.Range("D20")
$d20.NumberFormat = "@"
$d20.Value = "10.88"
$d20.Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "

$d21 = $ws.Range("D21")
$d21.NumberFormat = "@"
$d21.Value = "190.02"
$d21.Style = "Normal"
$ws.Range("E21").Value = "  +0.63%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("E23").Value = "  -0.01%  "

$d24 = $ws.Range("D24")
$d24.NumberFormat = "@"
$d24.Value = "145.18"
$d24.Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "

$d25 = $ws.Range("D25")
$d25.NumberFormat = "@"
$d25.Value = "7.782"
$d25.Style = "Normal"
$ws.Range("E25").Value = "  +1.93%  "

$d26 = $ws.Range("D26")
$d26.NumberFormat = "@"
$d26.Value = "0.1247"
$d26.Style = "Normal"
$ws.Range("E26").Value = "  +6.04%  "

$d27 = $ws.Range("D27")
$d27.NumberFormat = "@"
$d27.Value = "15.84"
$d27.Style = "Normal"
$ws.Range("E27").Value = "  +1.42%  "

$d28 = $ws.Range("D28")
$d28.NumberFormat = "@"
$d28.Value = "0.06646"
$d28.Style = "Normal"
$ws.Range("E28").Value = "  +11.72%  "

$d29 = $ws.Range("D29")
$d29.NumberFormat = "@"
$d29.Value = "1.358"
$d29.Style = "Normal"
$ws.Range("E29").Value = "  +5.63%  "

$ws.Range("E30").Value = "  +0.42%  "

$d31 = $ws.Range("D31")
$d31.NumberFormat = "@"
$d31.Value = "3.593"
$d31.Style = "Normal"
$ws.Range("E31").Value = "  +2.08%  "

$d32 = $ws.Range("D32")
$d32.NumberFormat = "@"
$d32.Value = "3.579"
$d32.Style = "Normal"
$ws.Range("E32").Value = "  +1.68%  "

$d33 = $ws.Range("D33")
$d33.NumberFormat = "@"
$d33.Value = "1.663"
$d33.Style = "Normal"
$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("E34").Value = "  +1.54%  "

$ws.Range("E35").Value = "  +3.25%  "

$ws.Range("E36").Value = "  +0.81%  "

$d37 = $ws.Range("D37")
$d37.NumberFormat = "@"
$d37.Value = "2.726"
$d37.Style = "Normal"
$ws.Range("E37").Value = "  +2.72%  "

$d38 = $ws.Range("D38")
$d38.NumberFormat = "@"
$d38.Value = "6.444"
$d38.Style = "Normal"
$ws.Range("E38").Value = "  +5.81%  "

$ws.Range("E39").Value = "  +0.58%  "

$d40 = $ws.Range("D40")
$d40.NumberFormat = "@"
$d40.Value = "1.108.58"
$d40.Style = "Normal"
$ws.Range("E40").Value = "  +2.93%  "

$d41 = $ws.Range("D41")
$d41.NumberFormat = "@"
$d41.Value = "0.8799"
$d41.Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "

$ws.Range("E42").Value = "  +0.45%  "

$d43 = $ws.Range("D43")
$d43.NumberFormat = "@"
$d43.Value = "100.81"
$d43.Style = "Normal"
$ws.Range("E43").Value = "  +0.66%  "

$d44 = $ws.Range("D44")
$d44.NumberFormat = "@"
$d44.Value = "1.835.40"
$d44.Style = "Normal"
$ws.Range("E44").Value = "  +0.81%  "

$d45 = $ws.Range("D45")
$d45.NumberFormat = "@"
$d45.Value = "0.00000000114"
$d45.Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

$d46 = $ws.Range("D46")
$d46.NumberFormat = "@"
$d46.Value = "57.01"
$d46.Style = "Normal"
$ws.Range("E46").Value = "  +1.47%  "

$d47 = $ws.Range("D47")
$d47.NumberFormat = "@"
$d47.Value = "1.010"
$d47.Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "

$d48 = $ws.Range("D48")
$d48.NumberFormat = "@"
$d48.Value = "8.137"
$d48.Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "

$d49 = $ws.Range("D49")
$d49.NumberFormat = "@"
$d49.Value = "0.05273"
$d49.Style = "Normal"
$ws.Range("E49").Value = "  +1.22%  "

$d50 = $ws.Range("D50")
$d50.NumberFormat = "@"
$d50.Value = "0.4299"
$d50.Style = "Normal"

$d51 = $ws.Range("D51")
$d51.NumberFormat = "@"
$d51.Value = "6.074"
$d51.Style = "Normal"
$ws.Range("E51").Value = "  +3.37%  "

